$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.756.30'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.744.47'
$ws.Range('D3').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.41'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5035'
$ws.Range('D7').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2716'
$ws.Range('D9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06161'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.745.60'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06937'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.50'
$ws.Range('D13').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6027'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.06'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.761.42'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006871'
$ws.Range('D20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.59'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.966.22'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.049'
$ws.Range('D23').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.161'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.76'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.481'
$ws.Range('D27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.818'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.03'
$ws.Range('D29').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08123'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.741'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.491'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04564'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9993'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.618'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9845'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6113'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.670'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01557'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.917'
$ws.Range('D41').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.91'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3839'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.073'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7323'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05372'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1114'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.976'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.23'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.53'
$ws.Range('D51').ClearFormats()

$ws.Range('E2').Value = '  -2.67%  '
$ws.Range('E3').Value = '  -4.98%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -8.81%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  -6.39%  '
$ws.Range('E8').Value = '  -6.41%  '
$ws.Range('E9').Value = '  -9.47%  '
$ws.Range('E10').Value = '  -10.62%  '
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('E13').Value = '  -11.98%  '
$ws.Range('E14').Value = '  -9.51%  '
$ws.Range('E15').Value = '  -18.04%  '
$ws.Range('E16').Value = '  -13.35%  '
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('E20').Value = '  -13.14%  '
$ws.Range('E21').Value = '  -16.32%  '
$ws.Range('E22').Value = '  -5.44%  '
$ws.Range('E23').Value = '  -11.92%  '
$ws.Range('E25').Value = '  -11.33%  '
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('E27').Value = '  -13.93%  '
$ws.Range('E28').Value = '  -16.67%  '
$ws.Range('E29').Value = '  -11.56%  '
$ws.Range('E30').Value = '  -6.22%  '
$ws.Range('E31').Value = '  -8.11%  '
$ws.Range('E32').Value = '  -11.95%  '
$ws.Range('E33').Value = '  -13.62%  '
$ws.Range('E34').Value = '  -5.61%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -10.30%  '
$ws.Range('E37').Value = '  -12.99%  '
$ws.Range('E38').Value = '  -16.21%  '
$ws.Range('E39').Value = '  -13.56%  '
$ws.Range('E40').Value = '  -9.18%  '
$ws.Range('E41').Value = '  -14.75%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  -5.71%  '
$ws.Range('E44').Value = '  -18.53%  '
$ws.Range('E45').Value = '  -13.94%  '
$ws.Range('E46').Value = '  -18.84%  '
$ws.Range('E47').Value = '  -6.87%  '
$ws.Range('E48').Value = '  -10.93%  '
$ws.Range('E49').Value = '  -19.17%  '
$ws.Range('E50').Value = '  -13.15%  '
$ws.Range('E51').Value = '  -12.54%  '
